# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

# Update the "Conversión del día" note with the new Binance rates.
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.67 = 22930.26 pesos`n✅ 22930.26 pesos = 5.65 = 949.64 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Update the underlying rate table on the "tasas" sheet.
$wsTasas.Range("N10").Value = 176.23
$wsTasas.Range("O10").Value = 4041
$wsTasas.Range("N12").Value = 4059.5
$wsTasas.Range("O12").Value = 168.121
